# Rename the paired-column headers so the "_old"/"_new" suffixes reflect the
# actual format-version identifiers being compared (FV2404 -> FV2410).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404",
  "Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410",
  "Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# Turn the data range into a real Excel Table (ListObject) so it can be
# filtered/sorted consistently, matching the column headers above.
$tableRange = $ws.Range("A1:U82")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row so it stays visible while scrolling.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
